# Scheduled-runner update: refresh cached market-price / profit figures
# (currentAveragePrice, currentAveragePriceNQ, LevePriceNQ, LeveProfitNQ, etc.)
# across the per-job-sheet profit tables (Table_<JOB>) in this workbook.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H70").Value = 18326.834
$ws.Range("I70").Value = 1991.6
$ws.Range("K70").Value = 5974.799999999999
$ws.Range("M70").Value = -5704.799999999999
$ws.Range("H73").Value = 18326.834
$ws.Range("I73").Value = 1991.6
$ws.Range("K73").Value = 5974.799999999999
$ws.Range("M73").Value = -5038.799999999999
$ws.Range("H129").Value = 2163.8
$ws.Range("I129").Value = 1848.6666
$ws.Range("K129").Value = 5545.9998
$ws.Range("M129").Value = -545.9997999999996
$ws.Range("H132").Value = 45060.945
$ws.Range("I132").Value = 47716.57
$ws.Range("K132").Value = 143149.71
$ws.Range("M132").Value = -140619.71
$ws.Range("H137").Value = 1427256.6
$ws.Range("I137").Value = 1390321.1
$ws.Range("K137").Value = 4170963.3
$ws.Range("M137").Value = -4168413.3
$ws.Range("H138").Value = 4574.6772
$ws.Range("I138").Value = 6751.75
$ws.Range("J138").Value = 3817.4348
$ws.Range("K138").Value = 20255.25
$ws.Range("L138").Value = 11452.3044
$ws.Range("M138").Value = -15115.25
$ws.Range("N138").Value = -21732.3044

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1454.3726
$ws.Range("I2").Value = 1316.2683
$ws.Range("J2").Value = 2020.6
$ws.Range("K2").Value = 1316.2683
$ws.Range("L2").Value = 2020.6
$ws.Range("M2").Value = -1203.2683
$ws.Range("N2").Value = -2246.6
$ws.Range("H32").Value = 2119.52
$ws.Range("I32").Value = 2066.5154
$ws.Range("K32").Value = 2066.5154
$ws.Range("M32").Value = -1779.5154
$ws.Range("H61").Value = 702591.2
$ws.Range("I61").Value = 986599.0600000001
$ws.Range("K61").Value = 986599.0600000001
$ws.Range("M61").Value = -986387.0600000001
$ws.Range("H63").Value = 3796.6924
$ws.Range("I63").Value = 2040
$ws.Range("K63").Value = 2040
$ws.Range("M63").Value = -1354
$ws.Range("H66").Value = 3796.6924
$ws.Range("I66").Value = 2040
$ws.Range("K66").Value = 10200
$ws.Range("M66").Value = -6768
$ws.Range("H74").Value = 2275450.2
$ws.Range("I74").Value = 2606554.8
$ws.Range("K74").Value = 2606554.8
$ws.Range("M74").Value = -2605680.8
$ws.Range("H77").Value = 2275450.2
$ws.Range("I77").Value = 2606554.8
$ws.Range("K77").Value = 13032774
$ws.Range("M77").Value = -13028406
$ws.Range("H116").Value = 1454.3726
$ws.Range("I116").Value = 1316.2683
$ws.Range("J116").Value = 2020.6
$ws.Range("K116").Value = 1316.2683
$ws.Range("L116").Value = 2020.6
$ws.Range("M116").Value = 977.7317
$ws.Range("N116").Value = -6608.6
$ws.Range("H122").Value = 901.3182
$ws.Range("I122").Value = 901.3182
$ws.Range("K122").Value = 2703.9546
$ws.Range("M122").Value = -253.9546
$ws.Range("H132").Value = 281048.75
$ws.Range("I132").Value = 467651.97
$ws.Range("K132").Value = 1402955.91
$ws.Range("M132").Value = -1400425.91
$ws.Range("H136").Value = 702591.2
$ws.Range("I136").Value = 986599.0600000001
$ws.Range("K136").Value = 2959797.18
$ws.Range("M136").Value = -2957247.18

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1454.3726
$ws.Range("I3").Value = 1316.2683
$ws.Range("J3").Value = 2020.6
$ws.Range("K3").Value = 1316.2683
$ws.Range("L3").Value = 2020.6
$ws.Range("M3").Value = -1202.2683
$ws.Range("N3").Value = -2248.6
$ws.Range("H86").Value = 2880.4666
$ws.Range("I86").Value = 1356.3334
$ws.Range("K86").Value = 1356.3334
$ws.Range("M86").Value = -233.3334
$ws.Range("H89").Value = 2880.4666
$ws.Range("I89").Value = 1356.3334
$ws.Range("K89").Value = 6781.666999999999
$ws.Range("M89").Value = -1165.666999999999
$ws.Range("H94").Value = 666.0417
$ws.Range("I94").Value = 651.5217
$ws.Range("K94").Value = 651.5217
$ws.Range("M94").Value = -200.5217
$ws.Range("H99").Value = 2557.4167
$ws.Range("I99").Value = 966.25
$ws.Range("J99").Value = 5739.75
$ws.Range("K99").Value = 966.25
$ws.Range("L99").Value = 5739.75
$ws.Range("M99").Value = 531.75
$ws.Range("N99").Value = -8735.75
$ws.Range("H107").Value = 1512.5769
$ws.Range("I107").Value = 1514.5217
$ws.Range("K107").Value = 1514.5217
$ws.Range("M107").Value = 405.4783
$ws.Range("H134").Value = 436848.9
$ws.Range("I134").Value = 796738
$ws.Range("K134").Value = 2390214
$ws.Range("M134").Value = -2387679

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 6752
$ws.Range("I31").Value = 1569.1666
$ws.Range("K31").Value = 1569.1666
$ws.Range("M31").Value = -1274.1666
$ws.Range("H34").Value = 6752
$ws.Range("I34").Value = 1569.1666
$ws.Range("K34").Value = 1569.1666
$ws.Range("M34").Value = -1367.1666
$ws.Range("H62").Value = 4000
$ws.Range("I62").Value = 3000
$ws.Range("K62").Value = 3000
$ws.Range("M62").Value = -2376
$ws.Range("H65").Value = 4000
$ws.Range("I65").Value = 3000
$ws.Range("K65").Value = 15000
$ws.Range("M65").Value = -11880
$ws.Range("H132").Value = 11923640
$ws.Range("I132").Value = 25336.732
$ws.Range("J132").Value = 41669400
$ws.Range("K132").Value = 76010.196
$ws.Range("L132").Value = 125008200
$ws.Range("M132").Value = -73480.196
$ws.Range("N132").Value = -125013260
$ws.Range("H134").Value = 3759.75
$ws.Range("I134").Value = 2462
$ws.Range("K134").Value = 7386
$ws.Range("M134").Value = -4851

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H52").Value = 1426.3334
$ws.Range("J52").Value = 1426.3334
$ws.Range("L52").Value = 4279.0002
$ws.Range("N52").Value = -4811.0002
$ws.Range("H68").Value = 6950.6313
$ws.Range("I68").Value = 2571
$ws.Range("J68").Value = 8118.533
$ws.Range("K68").Value = 7713
$ws.Range("L68").Value = 24355.599
$ws.Range("M68").Value = -6902
$ws.Range("N68").Value = -25977.599
$ws.Range("H71").Value = 6950.6313
$ws.Range("I71").Value = 2571
$ws.Range("J71").Value = 8118.533
$ws.Range("K71").Value = 23139
$ws.Range("L71").Value = 73066.79700000001
$ws.Range("M71").Value = -19083
$ws.Range("N71").Value = -81178.79700000001
$ws.Range("H131").Value = 16898.938
$ws.Range("J131").Value = 16692.2
$ws.Range("L131").Value = 50076.60000000001
$ws.Range("N131").Value = -60156.60000000001
$ws.Range("H132").Value = 2297.0256
$ws.Range("I132").Value = 2072.25
$ws.Range("J132").Value = 2396.926
$ws.Range("K132").Value = 18650.25
$ws.Range("L132").Value = 21572.334
$ws.Range("M132").Value = -16120.25
$ws.Range("N132").Value = -26632.334
$ws.Range("H137").Value = 3080.8064
$ws.Range("I137").Value = 1660.6842
$ws.Range("K137").Value = 4982.0526
$ws.Range("M137").Value = 117.9474
$ws.Range("H139").Value = 971.5333000000001
$ws.Range("I139").Value = 756.1667
$ws.Range("K139").Value = 2268.5001
$ws.Range("M139").Value = 2871.4999

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 458391.8
$ws.Range("I80").Value = 503471
$ws.Range("K80").Value = 503471
$ws.Range("M80").Value = -502473
$ws.Range("H83").Value = 458391.8
$ws.Range("I83").Value = 503471
$ws.Range("K83").Value = 2517355
$ws.Range("M83").Value = -2512363
$ws.Range("H102").Value = 2442.9707
$ws.Range("J102").Value = 4359.3335
$ws.Range("L102").Value = 4359.3335
$ws.Range("N102").Value = -7603.3335
$ws.Range("H126").Value = 837443.1
$ws.Range("I126").Value = 1391206.2
$ws.Range("J126").Value = 6798.375
$ws.Range("K126").Value = 4173618.6
$ws.Range("L126").Value = 20395.125
$ws.Range("M126").Value = -4171148.6
$ws.Range("N126").Value = -25335.125
$ws.Range("H132").Value = 189737.66
$ws.Range("I132").Value = 275234.34
$ws.Range("J132").Value = 1644.95
$ws.Range("K132").Value = 825703.02
$ws.Range("L132").Value = 4934.85
$ws.Range("M132").Value = -823173.02
$ws.Range("N132").Value = -9994.85

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 631048.25
$ws.Range("I132").Value = 737986.25
$ws.Range("K132").Value = 2213958.75
$ws.Range("M132").Value = -2211428.75
$ws.Range("H136").Value = 5442.3955
$ws.Range("I136").Value = 5020.9062
$ws.Range("K136").Value = 15062.7186
$ws.Range("M136").Value = -12512.7186

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 4666.3335
$ws.Range("I62").Value = 4500
$ws.Range("K62").Value = 4500
$ws.Range("M62").Value = -3876
$ws.Range("H65").Value = 4666.3335
$ws.Range("I65").Value = 4500
$ws.Range("K65").Value = 22500
$ws.Range("M65").Value = -19380
$ws.Range("H136").Value = 8088926.5
$ws.Range("J136").Value = 2923.818
$ws.Range("L136").Value = 8771.454000000002
$ws.Range("N136").Value = -13871.454
